# Add new PBAC permission - "When config changed" - to the Data Submission
# Emails section, and shift the following "Miscellaneous" rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 36; existing rows 36-39 shift down to 37-40.
$ws.Rows("36:36").Insert()

# Bring over the cell formatting (fills/borders) from the row that is now
# directly below (row 37, the old row 36) so the new row matches the look
# of the rest of the table.
$ws.Range("A37:X37").Copy()
$ws.Range("A36:X36").PasteSpecial(-4122)

# Fill in the values for the new permission row. Group stays "Data
# Submission Emails" (same as the rows above it).
$ws.Range("A36").Value = "Data Submission Emails"
$ws.Range("C36").Value = "Configuration is changed"
$ws.Range("X36").Value = "data_submission:cfg_changed"
$ws.Range("B36").Value = "When config changed"

$ws.Range("D36").Value = "unchecked"
$ws.Range("E36").Value = "checked"
$ws.Range("F36").Value = "unchecked"
$ws.Range("G36").Value = "checked"
$ws.Range("H36").Value = "fixed_unchecked"
$ws.Range("I36").Value = "fixed_unchecked"
$ws.Range("J36").Value = "fixed_checked"
$ws.Range("K36").Value = "fixed_unchecked"
$ws.Range("L36").Value = "fixed_checked"
$ws.Range("M36").Value = "fixed_unchecked"
$ws.Range("N36").Value = "checked"
$ws.Range("O36").Value = "checked"
$ws.Range("P36").Value = "checked"
$ws.Range("Q36").Value = "checked"
$ws.Range("R36").Value = "fixed_unchecked"
$ws.Range("S36").Value = "unchecked"
$ws.Range("T36").Value = "unchecked"
$ws.Range("U36").Value = "unchecked"
$ws.Range("V36").Value = "unchecked"
$ws.Range("W36").Value = "fixed_unchecked"

# Match the final saved selection/view seen in the authored workbook.
$ws.Range("V36").Select() | Out-Null
